$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Value)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "64.059.54"
Set-TextValue $ws "E2" "  +0.11%  "
Set-TextValue $ws "D3" "2.642.16"
Set-TextValue $ws "E3" "  +0.21%  "
Set-TextValue $ws "E4" "  -0.01%  "
Set-TextValue $ws "D5" "580.73"
Set-TextValue $ws "E5" "  +0.20%  "
Set-TextValue $ws "D6" "156.95"
Set-TextValue $ws "E6" "  +0.53%  "
Set-TextValue $ws "D7" "0.629"
Set-TextValue $ws "E7" "  -3.43%  "
Set-TextValue $ws "E8" "  +0.01%  "
Set-TextValue $ws "D9" "2.639.94"
Set-TextValue $ws "E9" "  +0.19%  "
Set-TextValue $ws "E10" "  -2.77%  "
Set-TextValue $ws "D11" "5.83"
Set-TextValue $ws "E11" "  +0.31%  "
Set-TextValue $ws "D12" "0.384"
Set-TextValue $ws "E12" "  -0.93%  "
Set-TextValue $ws "E13" "  +0.80%  "
Set-TextValue $ws "E14" "  +0.22%  "
Set-TextValue $ws "D15" "3.118.29"
Set-TextValue $ws "E15" "  +0.11%  "
Set-TextValue $ws "E16" "  +0.09%  "
Set-TextValue $ws "D17" "63.958.10"
Set-TextValue $ws "E17" "  +0.12%  "
Set-TextValue $ws "D18" "2.646.46"
Set-TextValue $ws "E18" "  -0.80%  "
Set-TextValue $ws "D19" "12.18"
Set-TextValue $ws "E19" "  -0.41%  "
Set-TextValue $ws "D20" "7.78"
Set-TextValue $ws "E20" "  +2.48%  "
Set-TextValue $ws "D21" "4.52"
Set-TextValue $ws "E21" "  -3.24%  "
Set-TextValue $ws "D22" "345.83"
Set-TextValue $ws "E22" "  -0.60%  "
Set-TextValue $ws "E23" "  +0.14%  "
Set-TextValue $ws "D24" "68.29"
Set-TextValue $ws "E24" "  +0.75%  "
Set-TextValue $ws "D25" "1.88"
Set-TextValue $ws "E25" "  +7.32%  "
Set-TextValue $ws "E26" "  +3.28%  "
Set-TextValue $ws "D27" "9.34"
Set-TextValue $ws "E27" "  -0.54%  "
Set-TextValue $ws "E28" "  +3.48%  "
Set-TextValue $ws "D29" "583.85"
Set-TextValue $ws "E29" "  +0.90%  "
Set-TextValue $ws "D30" "8.20"
Set-TextValue $ws "E30" "  +2.86%  "
Set-TextValue $ws "E31" "  +0.14%  "
Set-TextValue $ws "D32" "1.00"
Set-TextValue $ws "E32" "  -0.19%  "
Set-TextValue $ws "D33" "2.06"
Set-TextValue $ws "E33" "  -1.12%  "
Set-TextValue $ws "E34" "  +1.00%  "
Set-TextValue $ws "D35" "6.64"
Set-TextValue $ws "E35" "  +2.07%  "
Set-TextValue $ws "E36" "  +3.07%  "
Set-TextValue $ws "D37" "0.403"
Set-TextValue $ws "E37" "  -2.19%  "
Set-TextValue $ws "E38" "  -1.15%  "
Set-TextValue $ws "E39" "  -0.01%  "
Set-TextValue $ws "E40" "  +1.49%  "
Set-TextValue $ws "D41" "153.42"
Set-TextValue $ws "E41" "  +0.73%  "
Set-TextValue $ws "E42" "  +7.28%  "
Set-TextValue $ws "D43" "0.999"
Set-TextValue $ws "E43" "  -0.02%  "
Set-TextValue $ws "B44" "Aave"
Set-TextValue $ws "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D44" "162.88"
Set-TextValue $ws "E44" "  +2.34%  "
Set-TextValue $ws "B45" "InjectiveProtocol"
Set-TextValue $ws "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D45" "24.21"
Set-TextValue $ws "E45" "  +4.16%  "
Set-TextValue $ws "B46" "Filecoin"
Set-TextValue $ws "C46" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D46" "3.92"
Set-TextValue $ws "E46" "  -1.96%  "
Set-TextValue $ws "B47" "Hedera"
Set-TextValue $ws "C47" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D47" "0.0590"
Set-TextValue $ws "E47" "  -1.41%  "
Set-TextValue $ws "B48" "Mantle"
Set-TextValue $ws "C48" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D48" "0.635"
Set-TextValue $ws "E48" "  +0.13%  "
Set-TextValue $ws "B49" "Stellar"
Set-TextValue $ws "C49" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D49" "0.100"
Set-TextValue $ws "E49" "  -2.70%  "
Set-TextValue $ws "B50" "VeChain"
Set-TextValue $ws "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D50" "0.0249"
Set-TextValue $ws "E50" "  -1.68%  "
Set-TextValue $ws "B51" "EnergySwap"
Set-TextValue $ws "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D51" "19.13"
Set-TextValue $ws "E51" "  -0.43%  "
